$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3.1
$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 2.74
$ws.Range("I2").Value = 3.5
$ws.Range("J2").Value = 2.58
$ws.Range("K2").Value = 2.98
$ws.Range("F3").Value = 4.7
$ws.Range("G3").Value = 6.6
$ws.Range("H3").Value = 1.57
$ws.Range("I3").Value = 1.71
$ws.Range("J3").Value = 4.6
$ws.Range("K3").Value = 5.6
$ws.Range("F4").Value = 3.25
$ws.Range("G4").Value = 3.6
$ws.Range("H4").Value = 2.02
$ws.Range("I4").Value = 2.2
$ws.Range("J4").Value = 4.1
$ws.Range("K4").Value = 4.7
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 1.16
$ws.Range("P4").Value = 2.72
$ws.Range("Q4").Value = 1.49
$ws.Range("R4").Value = 1.69
$ws.Range("S4").Value = 2.06
$ws.Range("T4").Value = 1.5
$ws.Range("U4").Value = 2.68
$ws.Range("X4").Value = 34
$ws.Range("Y4").Value = 18.5
$ws.Range("Z4").Value = 21
$ws.Range("AA4").Value = 30
$ws.Range("AB4").Value = 25
$ws.Range("AC4").Value = 12.5
$ws.Range("AD4").Value = 13.5
$ws.Range("AE4").Value = 22
$ws.Range("AF4").Value = 36
$ws.Range("AG4").Value = 18
$ws.Range("AH4").Value = 17.5
$ws.Range("AI4").Value = 29
$ws.Range("AJ4").Value = 70
$ws.Range("AK4").Value = 38
$ws.Range("AL4").Value = 40
$ws.Range("AM4").Value = 60
$ws.Range("AN4").Value = 23
$ws.Range("AO4").Value = 10.5
$ws.Range("F5").Value = 3.55
$ws.Range("G5").Value = 4.1
$ws.Range("H5").Value = 2.1
$ws.Range("I5").Value = 2.32
$ws.Range("K5").Value = 3.75
$ws.Range("P5").Value = 1.8
$ws.Range("Q5").Value = 2.04
$ws.Range("F6").Value = 1.83
$ws.Range("G6").Value = 2.04
$ws.Range("H6").Value = 4.1
$ws.Range("I6").Value = 7.8
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 3.6
$ws.Range("G8").Value = 7.8
$ws.Range("J8").Value = 4.3
$ws.Range("P8").Value = 2.2
$ws.Range("AI8").Value = 980
$ws.Range("F9").Value = 1.93
$ws.Range("H9").Value = 4.2
$ws.Range("I9").Value = 4.6
$ws.Range("J9").Value = 3.7
$ws.Range("P9").Value = 1.96
$ws.Range("Q9").Value = 1.78
$ws.Range("P10").Value = 2.22
$ws.Range("F11").Value = 5
$ws.Range("Q11").Value = 2.02
$ws.Range("G12").Value = 1.97
$ws.Range("Q12").Value = 1.61
$ws.Range("I13").Value = 14
$ws.Range("U13").Value = 1.72
$ws.Range("X13").Value = 24
$ws.Range("AJ13").Value = 10
$ws.Range("H14").Value = 8
$ws.Range("K14").Value = 5.2
$ws.Range("H15").Value = 1.79
$ws.Range("I15").Value = 1.86
$ws.Range("P15").Value = 2.08
$ws.Range("Q15").Value = 1.78
$ws.Range("F16").Value = 6.2
$ws.Range("Q16").Value = 1.76
$ws.Range("F17").Value = 1.46
$ws.Range("G17").Value = 1.52
$ws.Range("K17").Value = 5.6
$ws.Range("N17").Value = 5.4
$ws.Range("O17").Value = 1.19
$ws.Range("Q17").Value = 1.56
$ws.Range("R17").Value = 1.6
$ws.Range("U17").Value = 2.1
$ws.Range("AK17").Value = 15.5
$ws.Range("AN17").Value = 6.6
$ws.Range("J18").Value = 4
$ws.Range("G19").Value = 1.7
$ws.Range("U19").Value = 2.02
$ws.Range("X19").Value = 19
$ws.Range("Y19").Value = 20
$ws.Range("Z19").Value = 110
$ws.Range("AC19").Value = 10
$ws.Range("AE19").Value = 870
$ws.Range("AF19").Value = 9.800000000000001
$ws.Range("AG19").Value = 9.800000000000001
$ws.Range("H20").Value = 3.1
$ws.Range("O21").Value = 1.2
$ws.Range("P21").Value = 2.44
$ws.Range("Q21").Value = 1.62
$ws.Range("S21").Value = 2.5
$ws.Range("T21").Value = 1.96
$ws.Range("U21").Value = 1.93
$ws.Range("AF21").Value = 9
$ws.Range("AL21").Value = 80
$ws.Range("AN21").Value = 4.9
$ws.Range("F22").Value = 1.6
$ws.Range("I22").Value = 6.4
$ws.Range("J22").Value = 4.3
$ws.Range("P22").Value = 2.46
$ws.Range("G23").Value = 1.83
$ws.Range("H23").Value = 4.9
$ws.Range("P23").Value = 2.08
$ws.Range("H24").Value = 3.85
$ws.Range("I24").Value = 4.1
$ws.Range("J24").Value = 4.1
$ws.Range("Q24").Value = 1.72
$ws.Range("P25").Value = 2.66
$ws.Range("Q25").Value = 1.41
$ws.Range("K26").Value = 3.15
$ws.Range("P26").Value = 1.46
$ws.Range("Q26").Value = 2.84
$ws.Range("F27").Value = 2.38
$ws.Range("G27").Value = 2.64
$ws.Range("H27").Value = 3.2
$ws.Range("I27").Value = 3.7
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = 3.5
$ws.Range("P27").Value = 1.64
$ws.Range("F29").Value = 2.02
$ws.Range("G29").Value = 2.16
$ws.Range("H29").Value = 4.7
$ws.Range("I29").Value = 5.3
$ws.Range("J29").Value = 3
$ws.Range("K29").Value = 3.3
$ws.Range("P29").Value = 1.49
$ws.Range("Q29").Value = 2.46
$ws.Range("F30").Value = 2.12
$ws.Range("G30").Value = 2.28
$ws.Range("H30").Value = 4.4
$ws.Range("I30").Value = 5.1
$ws.Range("J30").Value = 2.92
$ws.Range("K30").Value = 3.2
$ws.Range("P30").Value = 1.44
$ws.Range("Q30").Value = 2.56
$ws.Range("G31").Value = 970
$ws.Range("N31").Value = 1.03
$ws.Range("O31").Value = 1.3
$ws.Range("Q31").Value = 1.3
$ws.Range("R31").Value = 1.12
$ws.Range("S31").Value = 1.3
$ws.Range("F32").Value = 1.86
$ws.Range("G32").Value = 2.08
$ws.Range("H32").Value = 3.9
$ws.Range("I32").Value = 6.4
$ws.Range("J32").Value = 3.2
$ws.Range("K32").Value = 3.9
$ws.Range("P32").Value = 1.5
$ws.Range("Q32").Value = 2.16
